$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format used by the data cells (border + explicit font) onto the
# header cell B1 so it matches the data column's formatting.
$ws.Range("B2").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Recolor the data values to an explicit black font.
$ws.Range("B2:B29").Font.Color = 0

# The data rows grow slightly taller (18.75 -> 19.5), the header row (1) stays put.
$ws.Range("A2:B29").RowHeight = 19.5
